$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of row 2 (date cell + plain "s=1" body cells)
# into the new row 7 by copying the whole row (values + styles); this keeps
# the existing cellXf/font/border indices instead of minting new style
# entries, matching the plain styling the new row should carry.
$ws.Range("A2:F2").Copy($ws.Range("A7:F7"))

$ws.Range("A7").Value = 46062
$ws.Range("B7").Value = "Creazione repository GitHub"
$ws.Range("C7").Value = "Creazione repository"
$ws.Range("D7").Value = "Gestione file"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = "Codice altervista ancora da caricare"

$tbl = $ws.ListObjects.Item("Table_1")
$tbl.Resize($ws.Range("A1:F7"))
